# KPI Exclusions Template - replace exclusion template content and fix sanity test
#
# The sheet held two near-duplicate "exclusion" rows (row 2 and row 3), each
# pairing a KPI name with a long semi-colon separated list of sub_categories
# to exclude. Row 3's list (the "linear_product_length_out_of_store" KPI) is
# the corrected/authoritative one; row 2's stale list is replaced with it
# (minus the erroneous "ALLERGY" entry), and the now-redundant row 3 is
# removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "COLD & FLU ; IBS ; KIDS COLD & FLU ; KIDS COUGH ;  KIDS DECONGESTANTS ; KIDS DIGESTIVE HEALTH ;  KIDS HAYFEVER ; LAXATIVES ; KIDS TEETHING ; WIND ;  REHYDRATION ; PROBITOTICS ; PAIN MANAGEMENT ; FIRST AID ; MIGRAINE RELIEF ; RASH TREATMENT ; ANTI-AGE FACE ; ARTIFICIAL TAN ; BODY CLEANSING ; COSMETICS ; DEODORANTS ; MEN'S TOILETRIES ; SUNCARE ; HAIR CARE ; BABY HEALTHCARE ; FOR MUM ; KIDS HAIRCARE ; KIDS TOILETRIES ; KIDS WIPES ; BABY SUNCARE ; COTTON ; INCONTINENCE ; SANITARY TOWELS ; FEMININE WASH ; KIDS MOUTHWASH"

$ws.Rows(3).Delete()

[void]$ws.Range("C2").Select()
